# Auto-generated edit script: updates market-data columns (H-N) on each
# per-job sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to reflect the
# latest scheduled-runner price pull. Cells with no new value (profit
# columns that no longer apply once price inputs changed) are cleared.

$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1891.0465
$ws.Range("I98").Value = 1518.0605
$ws.Range("K98").Value = 1518.0605
$ws.Range("M98").Value = -20.06050000000005
$ws.Range("H107").Value = 893.7222
$ws.Range("I107").Value = 761.5833
$ws.Range("J107").Value = 1158
$ws.Range("K107").Value = 761.5833
$ws.Range("L107").Value = 1158
$ws.Range("M107").Value = 1158.4167
$ws.Range("N107").Value = -4998
$ws.Range("H122").Value = 1891.0465
$ws.Range("I122").Value = 1518.0605
$ws.Range("K122").Value = 4554.181500000001
$ws.Range("M122").Value = -2104.181500000001
$ws.Range("H137").Value = 1716.7273
$ws.Range("J137").Value = 1842.1428
$ws.Range("L137").Value = 5526.428400000001
$ws.Range("N137").Value = -10626.4284
$ws.Range("H138").Value = 2374.7722
$ws.Range("J138").Value = 2217.889
$ws.Range("L138").Value = 6653.667
$ws.Range("N138").Value = -16933.667

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 191861.23
$ws.Range("I2").Value = 231727.17
$ws.Range("J2").Value = 504.8
$ws.Range("K2").Value = 231727.17
$ws.Range("L2").Value = 504.8
$ws.Range("M2").Value = -231614.17
$ws.Range("N2").Value = -730.8
$ws.Range("H32").Value = 3896.1125
$ws.Range("I32").Value = 2344.4722
$ws.Range("J32").Value = 17860.875
$ws.Range("K32").Value = 2344.4722
$ws.Range("L32").Value = 17860.875
$ws.Range("M32").Value = -2057.4722
$ws.Range("N32").Value = -18434.875
$ws.Range("H45").Value = 1591.6154
$ws.Range("I45").Value = 1249.75
$ws.Range("J45").Value = 1743.5555
$ws.Range("K45").Value = 1249.75
$ws.Range("L45").Value = 1743.5555
$ws.Range("M45").Value = -872.75
$ws.Range("N45").Value = -2497.5555
$ws.Range("H61").Value = 64411.07
$ws.Range("I61").Value = 89085.336
$ws.Range("J61").Value = 19997.4
$ws.Range("K61").Value = 89085.336
$ws.Range("L61").Value = 19997.4
$ws.Range("M61").Value = -88873.336
$ws.Range("N61").Value = -20421.4
$ws.Range("H74").Value = 822.61365
$ws.Range("I74").Value = 524.5854
$ws.Range("K74").Value = 524.5854
$ws.Range("M74").Value = 349.4146
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").Value = ""
$ws.Range("H77").Value = 822.61365
$ws.Range("I77").Value = 524.5854
$ws.Range("K77").Value = 2622.927
$ws.Range("M77").Value = 1745.073
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").Value = ""
$ws.Range("H110").Value = 1322
$ws.Range("I110").Value = 1322
$ws.Range("K110").Value = 1322
$ws.Range("M110").Value = 723
$ws.Range("H116").Value = 191861.23
$ws.Range("I116").Value = 231727.17
$ws.Range("J116").Value = 504.8
$ws.Range("K116").Value = 231727.17
$ws.Range("L116").Value = 504.8
$ws.Range("M116").Value = -229433.17
$ws.Range("N116").Value = -5092.8
$ws.Range("H122").Value = 1414.1666
$ws.Range("I122").Value = 1042.7273
$ws.Range("K122").Value = 3128.1819
$ws.Range("M122").Value = -678.1819
$ws.Range("H130").Value = 48907.582
$ws.Range("J130").Value = 48907.582
$ws.Range("L130").Value = 48907.582
$ws.Range("N130").Value = -58947.582
$ws.Range("H132").Value = 2327.96
$ws.Range("I132").Value = 1943.2273
$ws.Range("J132").Value = 5149.3335
$ws.Range("K132").Value = 5829.6819
$ws.Range("L132").Value = 15448.0005
$ws.Range("M132").Value = -3299.6819
$ws.Range("N132").Value = -20508.0005
$ws.Range("H136").Value = 64411.07
$ws.Range("I136").Value = 89085.336
$ws.Range("J136").Value = 19997.4
$ws.Range("K136").Value = 267256.008
$ws.Range("L136").Value = 59992.2
$ws.Range("M136").Value = -264706.008
$ws.Range("N136").Value = -65092.2

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 191861.23
$ws.Range("I3").Value = 231727.17
$ws.Range("J3").Value = 504.8
$ws.Range("K3").Value = 231727.17
$ws.Range("L3").Value = 504.8
$ws.Range("M3").Value = -231613.17
$ws.Range("N3").Value = -732.8
$ws.Range("H105").Value = 2204.182
$ws.Range("J105").Value = 3213
$ws.Range("L105").Value = 3213
$ws.Range("N105").Value = -6707
$ws.Range("H107").Value = 1510.6666
$ws.Range("I107").Value = 2017
$ws.Range("J107").Value = 1004.3333
$ws.Range("K107").Value = 2017
$ws.Range("L107").Value = 1004.3333
$ws.Range("M107").Value = -97
$ws.Range("N107").Value = -4844.3333
$ws.Range("H134").Value = 8659.666999999999
$ws.Range("I134").Value = 11283.077
$ws.Range("J134").Value = 1838.8
$ws.Range("K134").Value = 33849.231
$ws.Range("L134").Value = 5516.4
$ws.Range("M134").Value = -31314.231
$ws.Range("N134").Value = -10586.4

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 859.125
$ws.Range("I16").Value = 812.1667
$ws.Range("K16").Value = 812.1667
$ws.Range("M16").Value = -525.1667
$ws.Range("H31").Value = 2353.2222
$ws.Range("I31").Value = 1450.5555
$ws.Range("J31").Value = 3255.889
$ws.Range("K31").Value = 1450.5555
$ws.Range("L31").Value = 3255.889
$ws.Range("M31").Value = -1155.5555
$ws.Range("N31").Value = -3845.889
$ws.Range("H34").Value = 2353.2222
$ws.Range("I34").Value = 1450.5555
$ws.Range("J34").Value = 3255.889
$ws.Range("K34").Value = 1450.5555
$ws.Range("L34").Value = 3255.889
$ws.Range("M34").Value = -1248.5555
$ws.Range("N34").Value = -3659.889
$ws.Range("H113").Value = 859.125
$ws.Range("I113").Value = 812.1667
$ws.Range("K113").Value = 812.1667
$ws.Range("M113").Value = 1357.8333
$ws.Range("H122").Value = 3203.3
$ws.Range("I122").Value = 1668.1666
$ws.Range("K122").Value = 5004.4998
$ws.Range("M122").Value = -2554.4998
$ws.Range("H132").Value = 1449
$ws.Range("I132").Value = 1443.5834
$ws.Range("K132").Value = 4330.7502
$ws.Range("M132").Value = -1800.7502
$ws.Range("H134").Value = 5491.3335
$ws.Range("I134").Value = 4589.8
$ws.Range("J134").Value = 9999
$ws.Range("K134").Value = 13769.4
$ws.Range("L134").Value = 29997
$ws.Range("M134").Value = -11234.4
$ws.Range("N134").Value = -35067

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 19653.139
$ws.Range("J131").Value = 20111.547
$ws.Range("L131").Value = 60334.641
$ws.Range("N131").Value = -70414.641
$ws.Range("H132").Value = 12400.556
$ws.Range("I132").Value = 1015
$ws.Range("J132").Value = 52250
$ws.Range("K132").Value = 9135
$ws.Range("L132").Value = 470250
$ws.Range("M132").Value = -6605
$ws.Range("N132").Value = -475310

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").Value = ""
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").Value = ""
$ws.Range("H102").Value = 5521.75
$ws.Range("I102").Value = 5521.75
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 5521.75
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -3899.75
$ws.Range("N102").Value = ""
$ws.Range("H113").Value = 1471.7142
$ws.Range("I113").Value = 1143.7142
$ws.Range("K113").Value = 1143.7142
$ws.Range("M113").Value = 1026.2858
$ws.Range("H122").Value = 1259.0952
$ws.Range("I122").Value = 1202.4615
$ws.Range("K122").Value = 3607.3845
$ws.Range("M122").Value = -1157.3845
$ws.Range("H126").Value = 3145420.5
$ws.Range("I126").Value = 3971926.2
$ws.Range("K126").Value = 11915778.6
$ws.Range("M126").Value = -11913308.6
$ws.Range("H132").Value = 1133716.6
$ws.Range("I132").Value = 1426709.5
$ws.Range("K132").Value = 4280128.5
$ws.Range("M132").Value = -4277598.5

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4200.6
$ws.Range("I7").Value = 4000.75
$ws.Range("J7").Value = 5000
$ws.Range("K7").Value = 4000.75
$ws.Range("L7").Value = 5000
$ws.Range("M7").Value = -3888.75
$ws.Range("N7").Value = -5224
$ws.Range("H69").Value = 200000
$ws.Range("J69").Value = 200000
$ws.Range("L69").Value = 200000
$ws.Range("N69").Value = -201622
$ws.Range("H72").Value = 200000
$ws.Range("J72").Value = 200000
$ws.Range("L72").Value = 600000
$ws.Range("N72").Value = -608112
$ws.Range("H126").Value = 4200.6
$ws.Range("I126").Value = 4000.75
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 12002.25
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -9532.25
$ws.Range("N126").Value = -19940
$ws.Range("H132").Value = 2058.8298
$ws.Range("I132").Value = 1577.5416
$ws.Range("K132").Value = 4732.6248
$ws.Range("M132").Value = -2202.6248
$ws.Range("H136").Value = 2787.5
$ws.Range("I136").Value = 2244.6667
$ws.Range("K136").Value = 6734.000100000001
$ws.Range("M136").Value = -4184.000100000001

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 831.75
$ws.Range("I113").Value = 603.6667
$ws.Range("K113").Value = 1811.0001
$ws.Range("M113").Value = 358.9999
$ws.Range("H122").Value = 39776.523
$ws.Range("I122").Value = 63130.617
$ws.Range("K122").Value = 189391.851
$ws.Range("M122").Value = -186941.851
$ws.Range("H126").Value = 6267.304
$ws.Range("I126").Value = 7915.875
$ws.Range("K126").Value = 23747.625
$ws.Range("M126").Value = -21277.625
$ws.Range("H132").Value = 1205.9354
$ws.Range("I132").Value = 1203.75
$ws.Range("K132").Value = 3611.25
$ws.Range("M132").Value = -1081.25
$ws.Range("H136").Value = 26456728
$ws.Range("I136").Value = 39683450
$ws.Range("K136").Value = 119050350
$ws.Range("M136").Value = -119047800
